$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new column before column L ("P/l before tax") to hold the new
# "Exceptional items" metric reported for the quarter.
$ws.Columns("L:L").Insert()

# Header rows (row 1 = lowercase labels, row 2 = titlecase labels)
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"

# Match the bold/centered/bordered header style used by the rest of row 1/2
$ws.Range("L1").Style = $ws.Range("K1").Style
$ws.Range("L2").Style = $ws.Range("K2").Style

# New quarterly figures (most quarters have no exceptional items)
$ws.Range("L3").Value = -13.26
$ws.Range("L7").Value = -0.13
$ws.Range("L12").Value = -5
$ws.Range("L13").Value = -5
